$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.592.66"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.947.24"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "363.71"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.92"
$ws.Range("E6").Value = "  -3.86%  "
$ws.Range("E7").Value = "  -2.94%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -3.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.40"
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0845"
$ws.Range("E12").Value = "  -2.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.78"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.410.80"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("E15").Value = "  -4.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.954.45"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.980"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.542.89"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.33"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.25"
$ws.Range("E21").Value = "  -4.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.04"
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.73"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("E25").Value = "  -3.05%  "
$ws.Range("E26").Value = "  -5.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.43"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.35"
$ws.Range("E29").Value = "  -4.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.109"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.28"
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.09"
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("E33").Value = "  +5.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.41"
$ws.Range("E34").Value = "  -5.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.32"
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("E37").Value = "  -3.34%  "
$ws.Range("E38").Value = "  +4.12%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.23"
$ws.Range("E40").Value = "  -5.49%  "
$ws.Range("E41").Value = "  -4.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.28"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("E43").Value = "  -3.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.33"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.091.07"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("E47").Value = "  -6.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.30"
$ws.Range("E48").Value = "  -7.21%  "
$ws.Range("E49").Value = "  -4.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0317"
$ws.Range("E50").Value = "  -6.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.83"
$ws.Range("E51").Value = "  -3.05%  "
